$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People | Crab Oil (G=5484)
$ws.Range("H32").Value = 3580.95
$ws.Range("I32").Value = 3694.3333
$ws.Range("J32").Value = 3410.875
$ws.Range("K32").Value = 3694.3333
$ws.Range("L32").Value = 3410.875
$ws.Range("M32").Value = -3368.3333
$ws.Range("N32").Value = -4062.875

# Row 74: Adhesive of Antipathy | Wing Glue (G=5507)
$ws.Range("H74").Value = 12900.5
$ws.Range("J74").Value = 15667.333
$ws.Range("L74").Value = 15667.333
$ws.Range("N74").Value = -17539.333

# Row 77: It's Gonna Grow Back (L) | Wing Glue (G=5507)
$ws.Range("H77").Value = 12900.5
$ws.Range("J77").Value = 15667.333
$ws.Range("L77").Value = 78336.66500000001
$ws.Range("N77").Value = -87696.66500000001

# Row 80: Cleansing the Wicked Humours | Hallowed Water (G=12605)
$ws.Range("H80").Value = 3381
$ws.Range("I80").Value = 3549.6
$ws.Range("J80").Value = 3100
$ws.Range("K80").Value = 10648.8
$ws.Range("L80").Value = 9300
$ws.Range("M80").Value = -9650.799999999999
$ws.Range("N80").Value = -11296

# Row 83: Washing Away the Sins (L) | Hallowed Water (G=12605)
$ws.Range("H83").Value = 3381
$ws.Range("I83").Value = 3549.6
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 31946.4
$ws.Range("L83").Value = 27900
$ws.Range("M83").Value = -26954.4
$ws.Range("N83").Value = -37884

# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink (G=12603)
$ws.Range("H86").Value = 3227.6667
$ws.Range("I86").Value = 3469.8
$ws.Range("K86").Value = 3469.8
$ws.Range("M86").Value = -2346.8

# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink (G=12603)
$ws.Range("H89").Value = 3227.6667
$ws.Range("I89").Value = 3469.8
$ws.Range("K89").Value = 17349
$ws.Range("M89").Value = -11733

# Row 95: Official Strategy Guide | Gyuki Leather Codex (G=18200)
$ws.Range("H95").Value = 33634.5
$ws.Range("J95").Value = 33634.5
$ws.Range("L95").Value = 33634.5
$ws.Range("N95").Value = -39126.5

# Row 106: Making Your Mark | Enchanted Palladium Ink (G=19903)
$ws.Range("H106").Value = 16322.25
$ws.Range("I106").Value = 2861
$ws.Range("K106").Value = 2861
$ws.Range("M106").Value = -2230


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot (G=27713)
$ws.Range("H2").Value = 1240.6666
$ws.Range("I2").Value = 1202.6
$ws.Range("K2").Value = 1202.6
$ws.Range("M2").Value = -1089.6

# Row 32: Ingot We Trust | Steel Ingot (G=44147)
$ws.Range("H32").Value = 4292.882
$ws.Range("I32").Value = 3459.3176
$ws.Range("K32").Value = 3459.3176
$ws.Range("M32").Value = -3172.3176

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot (G=43999)
$ws.Range("H61").Value = 2472.1667
$ws.Range("I61").Value = 2432.0908
$ws.Range("K61").Value = 2432.0908
$ws.Range("M61").Value = -2220.0908

# Row 116: No Scope | Titanbronze Ingot (G=27713)
$ws.Range("H116").Value = 1240.6666
$ws.Range("I116").Value = 1202.6
$ws.Range("K116").Value = 1202.6
$ws.Range("M116").Value = 1091.4

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot (G=43997)
$ws.Range("H132").Value = 2790.8
$ws.Range("I132").Value = 2549.4075
$ws.Range("J132").Value = 4963.3335
$ws.Range("K132").Value = 7648.2225
$ws.Range("L132").Value = 14890.0005
$ws.Range("M132").Value = -5118.2225
$ws.Range("N132").Value = -19950.0005

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot (G=43999)
$ws.Range("H136").Value = 2472.1667
$ws.Range("I136").Value = 2432.0908
$ws.Range("K136").Value = 7296.2724
$ws.Range("M136").Value = -4746.2724


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot (G=27713)
$ws.Range("H3").Value = 1240.6666
$ws.Range("I3").Value = 1202.6
$ws.Range("K3").Value = 1202.6
$ws.Range("M3").Value = -1088.6

# Row 22: Riveting Run | Iron Rivets (G=5092)
$ws.Range("H22").Value = 1818.8
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents() | Out-Null

# Row 134: Ruthenium Supremium | Ruthenium Ingot (G=43998)
$ws.Range("H134").Value = 1623.5264
$ws.Range("I134").Value = 1619.5625
$ws.Range("J134").Value = 1644.6666
$ws.Range("K134").Value = 4858.6875
$ws.Range("L134").Value = 4933.9998
$ws.Range("M134").Value = -2323.6875
$ws.Range("N134").Value = -10003.9998


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 29: Grinding It Out | Mudstone Grinding Wheel (G=2408)
$ws.Range("H29").Value = 1668.125
$ws.Range("I29").Value = 1620.7142
$ws.Range("K29").Value = 1620.7142
$ws.Range("M29").Value = -1327.7142

# Row 132: Hull Lotta Damage | Ginseng Lumber (G=44019)
$ws.Range("H132").Value = 6605.7334
$ws.Range("I132").Value = 7553.4546
$ws.Range("K132").Value = 22660.3638
$ws.Range("M132").Value = -20130.3638

# Row 134: Wood You Be Quiet | Ceiba Lumber (G=44020)
$ws.Range("H134").Value = 2455.9636
$ws.Range("I134").Value = 2064.745
$ws.Range("K134").Value = 6194.235
$ws.Range("M134").Value = -3659.235


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up | Kukuru Butter (G=4854)
$ws.Range("H12").Value = 184.73334
$ws.Range("I12").Value = 65.25
$ws.Range("J12").Value = 228.18182
$ws.Range("K12").Value = 195.75
$ws.Range("L12").Value = 684.5454599999999
$ws.Range("M12").Value = -22.75
$ws.Range("N12").Value = -1030.54546

# Row 34: Fever Pitch | Chamomile Tea (G=4749)
$ws.Range("H34").Value = 874.6111
$ws.Range("J34").Value = 1186.3636
$ws.Range("L34").Value = 3559.0908
$ws.Range("N34").Value = -3727.0908

# Row 39: Bloody Good Tart, This | Blood Currant Tart (G=4712)
$ws.Range("H39").Value = 149033.14
$ws.Range("J39").Value = 7246.6
$ws.Range("L39").Value = 21739.8
$ws.Range("N39").Value = -22327.8

# Row 45: Don't Turn Up Your Nose | Sauerkraut (G=29501)
$ws.Range("H45").Value = 12999.833
$ws.Range("J45").Value = 12999.833
$ws.Range("L45").Value = 38999.499
$ws.Range("N45").Value = -40063.499

# Row 52: Made by Apple in Coerthas | Apple Juice (G=31902)
$ws.Range("H52").Value = 790
$ws.Range("J52").Value = 790
$ws.Range("L52").Value = 2370
$ws.Range("N52").Value = -2902

# Row 55: Pagan Pastries | Pastry Fish (G=4733)
$ws.Range("H55").Value = 10503749
$ws.Range("J55").Value = 31255798
$ws.Range("L55").Value = 93767394
$ws.Range("N55").Value = -93767748

# Row 62: Little Orphan Candy | Fig Bavarois (G=12845)
$ws.Range("H62").Value = 95115.09
$ws.Range("I62").Value = 500599.5
$ws.Range("J62").Value = 5007.4443
$ws.Range("K62").Value = 1501798.5
$ws.Range("L62").Value = 15022.3329
$ws.Range("M62").Value = -1501112.5
$ws.Range("N62").Value = -16394.3329

# Row 64: The Aroma of Faith | Baked Onion Soup (G=12861)
$ws.Range("H64").Value = 3008.0908
$ws.Range("J64").Value = 4285.5713
$ws.Range("L64").Value = 12856.7139
$ws.Range("N64").Value = -13396.7139

# Row 65: Confections of Confession (L) | Fig Bavarois (G=12845)
$ws.Range("H65").Value = 95115.09
$ws.Range("I65").Value = 500599.5
$ws.Range("J65").Value = 5007.4443
$ws.Range("K65").Value = 4505395.5
$ws.Range("L65").Value = 45066.9987
$ws.Range("M65").Value = -4501963.5
$ws.Range("N65").Value = -51930.9987

# Row 67: Soup's On (L) | Baked Onion Soup (G=12861)
$ws.Range("H67").Value = 3008.0908
$ws.Range("J67").Value = 4285.5713
$ws.Range("L67").Value = 12856.7139
$ws.Range("N67").Value = -14728.7139

# Row 108: Meet for Meat | Grilled Rail (G=27853)
$ws.Range("H108").Value = 1252.3334
$ws.Range("J108").Value = 1500
$ws.Range("L108").Value = 4500
$ws.Range("N108").Value = -10260

# Row 110: His Dark Utensils | Spaghetti al Nero (G=27857)
$ws.Range("H110").Value = 12814
$ws.Range("J110").Value = 16599.6
$ws.Range("L110").Value = 49798.8
$ws.Range("N110").Value = -57978.8

# Row 117: A Good Omen | Peppered Popotoes (G=27870)
$ws.Range("H117").Value = 3908.3447
$ws.Range("I117").Value = 1213.3334
$ws.Range("J117").Value = 5121.1
$ws.Range("K117").Value = 3640.0002
$ws.Range("L117").Value = 15363.3
$ws.Range("M117").Value = -198.0001999999999
$ws.Range("N117").Value = -22247.3


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57: Gold Is So Last Year | Electrum Circlet (Amber) (G=2876)
$ws.Range("H57").Value = 12110
$ws.Range("I57").Value = 12110
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 12110
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -11290
$ws.Range("N57").ClearContents() | Out-Null

# Row 113: Copious Crystal Cannons | Manasilver Nugget (G=27710)
$ws.Range("H113").Value = 2854.3333
$ws.Range("I113").Value = 2901.6667
$ws.Range("K113").Value = 2901.6667
$ws.Range("M113").Value = -731.6667000000002

# Row 132: On Board for Lar | Lar Ingot (G=44008)
$ws.Range("H132").Value = 4483.025
$ws.Range("I132").Value = 2802.3667
$ws.Range("K132").Value = 8407.1001
$ws.Range("M132").Value = -5877.1001


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather (G=5277)
$ws.Range("H22").Value = 3729.3215
$ws.Range("J22").Value = 2099.6667
$ws.Range("L22").Value = 2099.6667
$ws.Range("N22").Value = -2689.6667

# Row 27: Fire and Hide | Aldgoat Leather (G=5277)
$ws.Range("H27").Value = 3729.3215
$ws.Range("J27").Value = 2099.6667
$ws.Range("L27").Value = 2099.6667
$ws.Range("N27").Value = -2313.6667

# Row 40: Best Served Toad | Toad Leather (G=36248)
$ws.Range("H40").Value = 14318.134
$ws.Range("I40").Value = 7067.1304
$ws.Range("J40").Value = 38142.855
$ws.Range("K40").Value = 7067.1304
$ws.Range("L40").Value = 38142.855
$ws.Range("M40").Value = -6931.1304
$ws.Range("N40").Value = -38414.855

# Row 95: Weathering Heights | Gyuki Leather Highboots of Striking (G=18221)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents() | Out-Null

# Row 136: Respect for Br'aax | Br'aax Leather (G=44060)
$ws.Range("H136").Value = 4774.3706
$ws.Range("I136").Value = 3630.5293
$ws.Range("K136").Value = 10891.5879
$ws.Range("M136").Value = -8341.5879


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 97: Getting a Leg Up | Ruby Cotton Gaskins of Striking (G=18220)
$ws.Range("H97").Value = 43500
$ws.Range("J97").Value = 43500
$ws.Range("L97").Value = 43500
$ws.Range("N97").Value = -45482

# Row 114: Hunting Season | Pixie Cotton Hat of Striking (G=25978)
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents() | Out-Null

# Row 132: Comfy Cabins | Snow Cotton Cloth (G=44029)
$ws.Range("H132").Value = 4530.4287
$ws.Range("I132").Value = 3600.9443
$ws.Range("J132").Value = 6203.5
$ws.Range("K132").Value = 10802.8329
$ws.Range("L132").Value = 18610.5
$ws.Range("M132").Value = -8272.832900000001
$ws.Range("N132").Value = -23670.5

